$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "getRaceIds" (row 2) as fully tested (was "Half" -> now "X")
$ws.Range("B2").Value = "X"

# Mark "removeRaceById" (row 5) as tested ("X")
$ws.Range("B5").Value = "X"
$ws.Range("B5").HorizontalAlignment = -4108  # xlCenter

# Mark "removeRaceByName" (row 33) as tested ("X")
$ws.Range("B33").Value = "X"
$ws.Range("B33").HorizontalAlignment = -4108  # xlCenter

# Update the active selection to match the final cursor position
$ws.Range("G20").Select()
